$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell holds plain text (prices/volumes formatted as strings in the source feed).
# Force text storage (NumberFormat "@") while writing so Excel does not auto-convert
# numeric-looking values (e.g. "11.00", "0.9840") into numbers and drop the trailing
# zero, then ClearFormats() to drop the temporary format override again so the cell
# style stays exactly as it was (no "s" attribute / default style).

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.874.30"
$c.ClearFormats()

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -0.94%  "
$c.ClearFormats()

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.879.22"
$c.ClearFormats()

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  -1.81%  "
$c.ClearFormats()

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  -0.38%  "
$c.ClearFormats()

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "324.54"
$c.ClearFormats()

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  -1.10%  "
$c.ClearFormats()

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4622"
$c.ClearFormats()

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  -0.81%  "
$c.ClearFormats()

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3881"
$c.ClearFormats()

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  -1.97%  "
$c.ClearFormats()

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07851"
$c.ClearFormats()

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  -2.01%  "
$c.ClearFormats()

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.9840"
$c.ClearFormats()

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  -2.76%  "
$c.ClearFormats()

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  -2.63%  "
$c.ClearFormats()

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.856.40"
$c.ClearFormats()

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -3.00%  "
$c.ClearFormats()

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "7.006"
$c.ClearFormats()

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  -2.16%  "
$c.ClearFormats()

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.675"
$c.ClearFormats()

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -2.25%  "
$c.ClearFormats()

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.06977"
$c.ClearFormats()

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  +0.07%  "
$c.ClearFormats()

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "88.59"
$c.ClearFormats()

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  -0.81%  "
$c.ClearFormats()

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -0.31%  "
$c.ClearFormats()

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000009931"
$c.ClearFormats()

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -2.57%  "
$c.ClearFormats()

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "16.96"
$c.ClearFormats()

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  -1.97%  "
$c.ClearFormats()

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  -0.40%  "
$c.ClearFormats()

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "28.866.68"
$c.ClearFormats()

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  -0.98%  "
$c.ClearFormats()

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.271"
$c.ClearFormats()

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  -2.73%  "
$c.ClearFormats()

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "11.00"
$c.ClearFormats()

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -1.48%  "
$c.ClearFormats()

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.102"
$c.ClearFormats()

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  +2.37%  "
$c.ClearFormats()

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "155.56"
$c.ClearFormats()

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -0.30%  "
$c.ClearFormats()

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "19.33"
$c.ClearFormats()

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  -1.89%  "
$c.ClearFormats()

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "5.910"
$c.ClearFormats()

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -0.83%  "
$c.ClearFormats()

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "117.86"
$c.ClearFormats()

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  -2.39%  "
$c.ClearFormats()

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  -6.58%  "
$c.ClearFormats()

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.9016"
$c.ClearFormats()

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -4.61%  "
$c.ClearFormats()

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "5.270"
$c.ClearFormats()

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -1.98%  "
$c.ClearFormats()

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.321"
$c.ClearFormats()

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  -2.28%  "
$c.ClearFormats()

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.248"
$c.ClearFormats()

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.05754"
$c.ClearFormats()

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  -2.23%  "
$c.ClearFormats()

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.169"
$c.ClearFormats()

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  -0.27%  "
$c.ClearFormats()

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02079"
$c.ClearFormats()

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  -1.29%  "
$c.ClearFormats()

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -0.50%  "
$c.ClearFormats()

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "7.634"
$c.ClearFormats()

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -6.09%  "
$c.ClearFormats()

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.5669"
$c.ClearFormats()

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  -3.01%  "
$c.ClearFormats()

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.1776"
$c.ClearFormats()

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "9.693"
$c.ClearFormats()

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  -3.93%  "
$c.ClearFormats()

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "11.99"
$c.ClearFormats()

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -0.47%  "
$c.ClearFormats()

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.214"
$c.ClearFormats()

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -3.50%  "
$c.ClearFormats()

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.5335"
$c.ClearFormats()

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  -2.62%  "
$c.ClearFormats()

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.07043"
$c.ClearFormats()

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  -2.23%  "
$c.ClearFormats()

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -3.62%  "
$c.ClearFormats()

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.546"
$c.ClearFormats()

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +1.73%  "
$c.ClearFormats()

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "112.34"
$c.ClearFormats()

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  -0.61%  "
$c.ClearFormats()

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.060"
$c.ClearFormats()

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -6.40%  "
$c.ClearFormats()

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "70.84"
$c.ClearFormats()

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  -1.71%  "
$c.ClearFormats()

